$wb = $excel.ActiveWorkbook

# ALC!row106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 10756619
$ws.Range("I106").Value = 15156053
$ws.Range("J106").Value = 2445.111
$ws.Range("K106").Value = 15156053
$ws.Range("L106").Value = 2445.111
$ws.Range("M106").Value = -15155422
$ws.Range("N106").Value = -3707.111

# ALC!row107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 10000376
$ws.Range("I107").Value = 10869887
$ws.Range("K107").Value = 10869887
$ws.Range("M107").Value = -10867967

# ALC!row131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 3286.6667
$ws.Range("I131").Value = 495
$ws.Range("J131").Value = 5520
$ws.Range("K131").Value = 1485
$ws.Range("L131").Value = 16560
$ws.Range("M131").Value = 3555
$ws.Range("N131").Value = -26640

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2792.33
$ws.Range("I32").Value = 2674.526
$ws.Range("J32").Value = 6601.3335
$ws.Range("K32").Value = 2674.526
$ws.Range("L32").Value = 6601.3335
$ws.Range("M32").Value = -2387.526
$ws.Range("N32").Value = -7175.3335

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1513.8511
$ws.Range("I74").Value = 1123.9678
$ws.Range("J74").Value = 2269.25
$ws.Range("K74").Value = 1123.9678
$ws.Range("L74").Value = 2269.25
$ws.Range("M74").Value = -249.9677999999999
$ws.Range("N74").Value = -4017.25

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1513.8511
$ws.Range("I77").Value = 1123.9678
$ws.Range("J77").Value = 2269.25
$ws.Range("K77").Value = 5619.839
$ws.Range("L77").Value = 11346.25
$ws.Range("M77").Value = -1251.839
$ws.Range("N77").Value = -20082.25

# ARM!row97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1252.95
$ws.Range("I97").Value = 1197.5294
$ws.Range("J97").Value = 1567
$ws.Range("K97").Value = 1197.5294
$ws.Range("L97").Value = 1567
$ws.Range("M97").Value = -701.5293999999999
$ws.Range("N97").Value = -2559

# ARM!row102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4116854
$ws.Range("I102").Value = 5292505
$ws.Range("K102").Value = 5292505
$ws.Range("M102").Value = -5290883

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1028326.1
$ws.Range("I122").Value = 1352410.1
$ws.Range("J122").Value = 2060
$ws.Range("K122").Value = 4057230.3
$ws.Range("L122").Value = 6180
$ws.Range("M122").Value = -4054780.3
$ws.Range("N122").Value = -11080

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1727087
$ws.Range("I132").Value = 2110.3
$ws.Range("J132").Value = 5560368.5
$ws.Range("K132").Value = 6330.900000000001
$ws.Range("L132").Value = 16681105.5
$ws.Range("M132").Value = -3800.900000000001
$ws.Range("N132").Value = -16686165.5

# ARM!row137
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# BSM!row86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2093.7144
$ws.Range("I86").Value = 1980.5
$ws.Range("J86").Value = 2376.75
$ws.Range("K86").Value = 1980.5
$ws.Range("L86").Value = 2376.75
$ws.Range("M86").Value = -857.5
$ws.Range("N86").Value = -4622.75

# BSM!row89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2093.7144
$ws.Range("I89").Value = 1980.5
$ws.Range("J89").Value = 2376.75
$ws.Range("K89").Value = 9902.5
$ws.Range("L89").Value = 11883.75
$ws.Range("M89").Value = -4286.5
$ws.Range("N89").Value = -23115.75

# BSM!row99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 83334620
$ws.Range("I99").Value = 111112080
$ws.Range("J99").Value = 2266.3333
$ws.Range("K99").Value = 111112080
$ws.Range("L99").Value = 2266.3333
$ws.Range("M99").Value = -111110582
$ws.Range("N99").Value = -5262.3333

# BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 10040.9
$ws.Range("I105").Value = 15401.5625
$ws.Range("J105").Value = 3914.4285
$ws.Range("K105").Value = 15401.5625
$ws.Range("L105").Value = 3914.4285
$ws.Range("M105").Value = -13654.5625
$ws.Range("N105").Value = -7408.4285

# BSM!row135
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 150000
$ws.Range("J135").Value = 150000
$ws.Range("L135").Value = 150000
$ws.Range("N135").Value = -160140

# BSM!row137
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 110000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# BSM!row138
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# BSM!row140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 47119.945
$ws.Range("J140").Value = 47119.945
$ws.Range("L140").Value = 47119.945
$ws.Range("N140").Value = -57479.945

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 374171.44
$ws.Range("I134").Value = 4054.8096
$ws.Range("J134").Value = 1669579.6
$ws.Range("K134").Value = 12164.4288
$ws.Range("L134").Value = 5008738.800000001
$ws.Range("M134").Value = -9629.4288
$ws.Range("N134").Value = -5013808.800000001

# CRP!row135
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# CUL!row5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1204.1282
$ws.Range("I5").Value = 409.7143
$ws.Range("J5").Value = 1649
$ws.Range("K5").Value = 1229.1429
$ws.Range("L5").Value = 4947
$ws.Range("M5").Value = -1117.1429
$ws.Range("N5").Value = -5171

# CUL!row109
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2215.3845
$ws.Range("I109").Value = 300
$ws.Range("J109").Value = 2375
$ws.Range("K109").Value = 900
$ws.Range("L109").Value = 7125
$ws.Range("M109").Value = 140
$ws.Range("N109").Value = -9205

# CUL!row122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 4071.0322
$ws.Range("I122").Value = 389.45
$ws.Range("J122").Value = 10764.818
$ws.Range("K122").Value = 3505.05
$ws.Range("L122").Value = 96883.36199999999
$ws.Range("M122").Value = -1055.05
$ws.Range("N122").Value = -101783.362

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2041743.4
$ws.Range("I131").Value = 6250544.5
$ws.Range("J131").Value = 1112.5151
$ws.Range("K131").Value = 18751633.5
$ws.Range("L131").Value = 3337.5453
$ws.Range("M131").Value = -18746593.5
$ws.Range("N131").Value = -13417.5453

# CUL!row135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1204.1282
$ws.Range("I135").Value = 409.7143
$ws.Range("J135").Value = 1649
$ws.Range("K135").Value = 3687.4287
$ws.Range("L135").Value = 14841
$ws.Range("M135").Value = -1152.4287
$ws.Range("N135").Value = -19911

# GSM!row113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 28572592
$ws.Range("I113").Value = 50000788
$ws.Range("J113").Value = 1664.4
$ws.Range("K113").Value = 50000788
$ws.Range("L113").Value = 1664.4
$ws.Range("M113").Value = -49998618
$ws.Range("N113").Value = -6004.4

# GSM!row135
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# GSM!row140
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 47578.684
$ws.Range("J140").Value = 47578.684
$ws.Range("L140").Value = 47578.684
$ws.Range("N140").Value = -57938.684

# LTW!row68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 111113150
$ws.Range("I68").Value = 2300
$ws.Range("J68").Value = 500001150
$ws.Range("K68").Value = 2300
$ws.Range("L68").Value = 500001150
$ws.Range("M68").Value = -1551
$ws.Range("N68").Value = -500002648

# LTW!row71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 111113150
$ws.Range("I71").Value = 2300
$ws.Range("J71").Value = 500001150
$ws.Range("K71").Value = 11500
$ws.Range("L71").Value = 2500005750
$ws.Range("M71").Value = -7756
$ws.Range("N71").Value = -2500013238

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10422433
$ws.Range("I132").Value = 11500134
$ws.Range("J132").Value = 4661.3335
$ws.Range("K132").Value = 34500402
$ws.Range("L132").Value = 13984.0005
$ws.Range("M132").Value = -34497872
$ws.Range("N132").Value = -19044.0005

# WVR!row107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 52632076
$ws.Range("I107").Value = 66667150
$ws.Range("J107").Value = 544
$ws.Range("K107").Value = 200001450
$ws.Range("L107").Value = 1632
$ws.Range("M107").Value = -199999530
$ws.Range("N107").Value = -5472

# WVR!row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1551.0588
$ws.Range("I122").Value = 1182.1666
$ws.Range("J122").Value = 2436.4
$ws.Range("K122").Value = 3546.4998
$ws.Range("L122").Value = 7309.200000000001
$ws.Range("M122").Value = -1096.4998
$ws.Range("N122").Value = -12209.2

# WVR!row135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 80715
$ws.Range("J135").Value = 80715
$ws.Range("L135").Value = 80715
$ws.Range("N135").Value = -90855

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2543.8
$ws.Range("I136").Value = 2621.4138
$ws.Range("J136").Value = 2403.125
$ws.Range("K136").Value = 7864.241399999999
$ws.Range("L136").Value = 7209.375
$ws.Range("M136").Value = -5314.241399999999
$ws.Range("N136").Value = -12309.375

# WVR!row137
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# WVR!row138
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# WVR!row139
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
